# D0 W data added
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CDF1")

for ($row = 2; $row -le 21; $row++) {
    $ws.Cells.Item($row, 2).Value = 80.379000000000005   # B: Q
    $ws.Cells.Item($row, 3).Value = $false                # C: Qmin
    $ws.Cells.Item($row, 4).Value = $false                # D: Qmax
}

# Update sheet view/selection to match final state
$ws.Range("G26").Select() | Out-Null
